$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The mapping table's row 3 (column "Kolumna") was renamed from
# "DelayGroup" to "DepDelayGroup" (the shared-string entry is dropped from
# its old slot and re-added at the end of the table, which is exactly what
# happens when Excel edits a cell's text in place).
$ws.Range("A3").Value = "DepDelayGroup"

# Reflect the author's final cursor position/selection on the sheet.
$ws.Range("A3").Select()
$excel.ActiveWindow.ScrollRow = 37
